$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.538.48'
$ws.Range("E2").Value = '  -0.18%  '
$ws.Range("D3").Value = '2.108.03'
$ws.Range("E3").Value = '  +10.14%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '251.69'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.22%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '47.84'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +8.67%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '59.39'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +2.12%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.374'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +1.91%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0747'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -1.91%  '
$ws.Range("E12").Value = '  +0.21%  '
$ws.Range("D13").Value = '2.411.48'
$ws.Range("E13").Value = '  +9.96%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '14.42'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.79%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.828'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +1.42%  '
$ws.Range("D16").Value = '2.109.40'
$ws.Range("E16").Value = '  +10.21%  '
$ws.Range("E17").Value = '  -0.09%  '
$ws.Range("D18").Value = '36.532.79'
$ws.Range("E18").Value = '  -0.07%  '
$ws.Range("E19").Value = '  -1.82%  '
$ws.Range("D20").Value = '0.0₃0832'
$ws.Range("E20").Value = '  -3.06%  '
$ws.Range("E21").Value = '  -0.31%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '240.22'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -3.86%  '
$ws.Range("E23").Value = '  -0.45%  '
$ws.Range("E24").Value = '  +0.16%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.45'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -7.56%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '171.87'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +2.47%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.47'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +14.82%  '
$ws.Range("E28").Value = '  +4.05%  '
$ws.Range("E29").Value = '  -9.60%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '28.45'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +62.04%  '
$ws.Range("E31").Value = '  -4.42%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.44'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -3.40%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0611'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -0.86%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0907'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +4.78%  '
$ws.Range("E35").Value = '  +15.05%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.945'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +7.29%  '
$ws.Range("E37").Value = '  -0.17%  '
$ws.Range("E38").Value = '  -4.04%  '
$ws.Range("E39").Value = '  -5.96%  '
$ws.Range("E40").Value = '  -12.00%  '
$ws.Range("E41").Value = '  +6.66%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0224'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -1.19%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '98.12'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -7.72%  '
$ws.Range("E44").Value = '  -1.94%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '16.16'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -7.04%  '
$ws.Range("D46").Value = '1.338.69'
$ws.Range("E46").Value = '  -0.01%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0843'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +3.53%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.01'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +9.90%  '
$ws.Range("D49").Value = '2.293.90'
$ws.Range("E49").Value = '  +9.88%  '
$ws.Range("E50").Value = '  +1.67%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.24'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -5.32%  '
